$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The cells contain a mojibake sequence (U+00C2 U+00B1, i.e. "Â±") where a
# correctly UTF-8 encoded "±" (U+00B1) character got double-decoded as
# Latin-1/CP1252. Build the search/replacement strings from explicit
# character codes so this works regardless of the script file's own encoding.
$mojibake = [string]([char]194) + [string]([char]177)
$fixed = [string]([char]177)

$rows = @(2,3,4,5,6,9,11,13,14,15,16,17)
$cols = @("B","C","D")

foreach ($r in $rows) {
    foreach ($col in $cols) {
        $addr = "$col$r"
        $cell = $ws.Range($addr)
        $val = $cell.Value2
        if ($val -ne $null -and $val.Contains($mojibake)) {
            $cell.Value2 = $val.Replace($mojibake, $fixed)
        }
    }
}
